# Apply the updated crypto price/volume snapshot to the sheet.
# Numeric-looking Price values are written via .Formula with a leading
# apostrophe (quote-prefix) so Excel stores them as text (preserving
# trailing zeros / exact formatting) instead of converting to a Double.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.267.90'
$ws.Range("E2").Value = '  +0.84%  '

$ws.Range("D3").Value = '3.428.32'
$ws.Range("E3").Value = '  +0.78%  '

$ws.Range("E4").Value = '  -0.03%  '

$ws.Range("D5").Formula = '''414.02'
$ws.Range("E5").Value = '  +1.37%  '

$ws.Range("D6").Formula = '''129.53'
$ws.Range("E6").Value = '  +1.39%  '

$ws.Range("E7").Value = '  -1.52%  '

$ws.Range("E8").Value = '  -0.01%  '

$ws.Range("E9").Value = '  -0.08%  '

$ws.Range("D10").Formula = '''0.140'
$ws.Range("E10").Value = '  +1.56%  '

$ws.Range("D11").Formula = '''42.91'
$ws.Range("E11").Value = '  +1.59%  '

$ws.Range("D12").Formula = '''9.20'
$ws.Range("E12").Value = '  +1.96%  '

$ws.Range("D13").Value = '3.973.15'
$ws.Range("E13").Value = '  +0.75%  '

$ws.Range("E14").Value = '  +6.27%  '

$ws.Range("E15").Value = '  -0.19%  '

$ws.Range("D16").Formula = '''20.46'
$ws.Range("E16").Value = '  -3.26%  '

$ws.Range("D17").Value = '3.457.05'
$ws.Range("E17").Value = '  +1.49%  '

$ws.Range("E18").Value = '  +4.74%  '

$ws.Range("E19").Value = '  +0.14%  '

$ws.Range("D20").Value = '62.289.27'
$ws.Range("E20").Value = '  +0.91%  '

$ws.Range("D21").Formula = '''467.79'
$ws.Range("E21").Value = '  +3.39%  '

$ws.Range("D22").Formula = '''91.13'
$ws.Range("E22").Value = '  -0.15%  '

$ws.Range("E23").Value = '  +3.63%  '

$ws.Range("D24").Formula = '''13.28'
$ws.Range("E24").Value = '  +3.53%  '

$ws.Range("B25").Value = 'Filecoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D25").Formula = '''10.28'
$ws.Range("E25").Value = '  +18.95%  '

$ws.Range("B26").Value = 'PancakeSwap'
$ws.Range("C26").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D26").Formula = '''3.32'
$ws.Range("E26").Value = '  +2.50%  '

$ws.Range("D27").Formula = '''33.01'
$ws.Range("E27").Value = '  -0.92%  '

$ws.Range("D28").Formula = '''4.78'
$ws.Range("E28").Value = '  +0.57%  '

$ws.Range("D29").Formula = '''7.76'
$ws.Range("E29").Value = '  +3.24%  '

$ws.Range("D30").Formula = '''11.88'
$ws.Range("E30").Value = '  -0.05%  '

$ws.Range("D31").Formula = '''2.64'
$ws.Range("E31").Value = '  -1.73%  '

$ws.Range("E32").Value = '  -0.42%  '

$ws.Range("E33").Value = '  -1.25%  '

$ws.Range("D34").Formula = '''41.15'
$ws.Range("E34").Value = '  -3.31%  '

$ws.Range("E35").Value = '  +0.05%  '

$ws.Range("D36").Formula = '''57.66'
$ws.Range("E36").Value = '  +8.58%  '

$ws.Range("E37").Value = '  -1.46%  '

$ws.Range("D38").Formula = '''0.999'
$ws.Range("E38").Value = '  -0.01%  '

$ws.Range("E39").Value = '  +5.28%  '

$ws.Range("D40").Formula = '''0.329'
$ws.Range("E40").Value = '  +5.21%  '

$ws.Range("E41").Value = '  +0.71%  '

$ws.Range("E42").Value = '  -0.20%  '

$ws.Range("D43").Formula = '''144.18'
$ws.Range("E43").Value = '  +2.55%  '

$ws.Range("B44").Value = 'ARBITRUM'
$ws.Range("C44").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D44").Formula = '''2.07'
$ws.Range("E44").Value = '  +5.78%  '

$ws.Range("B45").Value = 'WEMIXToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D45").Formula = '''2.64'
$ws.Range("E45").Value = '  +10.06%  '

$ws.Range("D46").Formula = '''4.34'
$ws.Range("E46").Value = '  +3.59%  '

$ws.Range("D47").Formula = '''2.43'
$ws.Range("E47").Value = '  +19.87%  '

$ws.Range("D48").Formula = '''16.48'
$ws.Range("E48").Value = '  +0.49%  '

$ws.Range("D49").Formula = '''22.20'
$ws.Range("E49").Value = '  -0.49%  '

$ws.Range("D50").Value = '0.0₃0515'
$ws.Range("E50").Value = '  +28.11%  '

$ws.Range("D51").Formula = '''111.41'
$ws.Range("E51").Value = '  +5.04%  '
